$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price cells keep their original text representation
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D33", "D34", "D37", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the data refresh
$ws.Range("D2").Value = '57.030.64'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '3.255.99'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '396.25'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '108.14'
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("E7").Value = '  +4.55%  '
$ws.Range("D9").Value = '0.619'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").Value = '39.35'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '0.0951'
$ws.Range("E11").Value = '  +5.98%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").Value = '3.763.73'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '8.29'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").Value = '18.95'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '3.266.78'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = '10.90'
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").Value = '56.856.34'
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("E21").Value = '  +5.53%  '
$ws.Range("D22").Value = '12.92'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").Value = '292.07'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = '74.16'
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").Value = '8.02'
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").Value = '28.13'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").Value = '0.169'
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").Value = '11.16'
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("D34").Value = '40.32'
$ws.Range("E34").Value = '  +10.47%  '
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("D37").Value = '51.24'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("E40").Value = '  -3.59%  '
$ws.Range("D41").Value = '137.57'
$ws.Range("E41").Value = '  +3.25%  '
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").Value = '3.92'
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '1.86'
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.282'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '16.59'
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").Value = '22.16'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D48").Value = '2.22'
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("D49").Value = '2.144.94'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '2.03'
$ws.Range("E50").Value = '  -4.71%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '2.32'
$ws.Range("E51").Value = '  -5.36%  '
